# Weekly update: insert a new observation as the most recent row (row 327)
# for "Vega Monumental Concepción" / Ajo, pushing the existing historical
# rows (327-355) down by one (to 328-356).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 327 (shifts 327:355 -> 328:356)
$ws.Rows.Item(327).Insert()

# Populate the new row with the latest weekly price observation
$ws.Range("A327").Value = 11
$ws.Range("B327").Value = "Vega Monumental Concepción"
$ws.Range("C327").Value = "Bíobío"
$ws.Range("D327").Value = 45265
$ws.Range("D327").NumberFormat = $ws.Range("D328").NumberFormat
$ws.Range("E327").Value = 8
$ws.Range("F327").Value = 100112003
$ws.Range("G327").Value = "Ajo"
$ws.Range("H327").Value = "Chino"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 300
$ws.Range("K327").Value = 22000
$ws.Range("L327").Value = 23000
$ws.Range("M327").Value = 22667
$ws.Range("N327").Value = "$/caja 10 kilos"
$ws.Range("O327").Value = "China"
$ws.Range("P327").Value = 2267
$ws.Range("Q327").Value = 10
$ws.Range("R327").Value = "Hortaliza"
